$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 640
$ws.Range("J18").Value = 675
$ws.Range("L18").Value = 675
$ws.Range("N18").Value = -1243
$ws.Range("H28").Value = 283.33334
$ws.Range("I28").Value = 348.33334
$ws.Range("J28").Value = 153.33333
$ws.Range("K28").Value = 348.33334
$ws.Range("L28").Value = 153.33333
$ws.Range("M28").Value = 136.66666
$ws.Range("N28").Value = -1123.33333
$ws.Range("H33").Value = 357.42856
$ws.Range("I33").Value = 310.8421
$ws.Range("K33").Value = 310.8421
$ws.Range("M33").Value = -81.84210000000002
$ws.Range("H92").Value = 627.3182
$ws.Range("I92").Value = 606.2353000000001
$ws.Range("K92").Value = 606.2353000000001
$ws.Range("M92").Value = 641.7646999999999
$ws.Range("H107").Value = 857.73914
$ws.Range("I107").Value = 898.5263
$ws.Range("K107").Value = 898.5263
$ws.Range("M107").Value = 1021.4737
$ws.Range("H116").Value = 4981.2
$ws.Range("J116").Value = 5312.4443
$ws.Range("L116").Value = 5312.4443
$ws.Range("N116").Value = -12196.4443
$ws.Range("H129").Value = 770.8570999999999
$ws.Range("I129").Value = 540
$ws.Range("J129").Value = 899.1111
$ws.Range("K129").Value = 1620
$ws.Range("L129").Value = 2697.3333
$ws.Range("M129").Value = 3380
$ws.Range("N129").Value = -12697.3333
$ws.Range("H137").Value = 17708.666
$ws.Range("I137").Value = 1962.3948
$ws.Range("J137").Value = 41643
$ws.Range("K137").Value = 5887.1844
$ws.Range("L137").Value = 124929
$ws.Range("M137").Value = -3337.1844
$ws.Range("N137").Value = -130029
$ws.Range("H138").Value = 2788.7886
$ws.Range("J138").Value = 2938
$ws.Range("L138").Value = 8814
$ws.Range("N138").Value = -19094

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3151.8845
$ws.Range("I45").Value = 2228.7778
$ws.Range("J45").Value = 3640.5881
$ws.Range("K45").Value = 2228.7778
$ws.Range("L45").Value = 3640.5881
$ws.Range("M45").Value = -1851.7778
$ws.Range("N45").Value = -4394.5881
$ws.Range("H74").Value = 3221.15
$ws.Range("I74").Value = 4205.909
$ws.Range("J74").Value = 2017.5555
$ws.Range("K74").Value = 4205.909
$ws.Range("L74").Value = 2017.5555
$ws.Range("M74").Value = -3331.909
$ws.Range("N74").Value = -3765.5555
$ws.Range("H77").Value = 3221.15
$ws.Range("I77").Value = 4205.909
$ws.Range("J77").Value = 2017.5555
$ws.Range("K77").Value = 21029.545
$ws.Range("L77").Value = 10087.7775
$ws.Range("M77").Value = -16661.545
$ws.Range("N77").Value = -18823.7775
$ws.Range("H132").Value = 18576.967
$ws.Range("I132").Value = 1616.6471
$ws.Range("K132").Value = 4849.9413
$ws.Range("M132").Value = -2319.9413

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 774.9167
$ws.Range("I107").Value = 754.4545000000001
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 754.4545000000001
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1165.5455
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 4389
$ws.Range("I134").Value = 5000.1333
$ws.Range("J134").Value = 1333.3334
$ws.Range("K134").Value = 15000.3999
$ws.Range("L134").Value = 4000.0002
$ws.Range("M134").Value = -12465.3999
$ws.Range("N134").Value = -9070.0002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10843.935
$ws.Range("I31").Value = 17748.75
$ws.Range("J31").Value = 3311.4092
$ws.Range("K31").Value = 17748.75
$ws.Range("L31").Value = 3311.4092
$ws.Range("M31").Value = -17453.75
$ws.Range("N31").Value = -3901.4092
$ws.Range("H34").Value = 10843.935
$ws.Range("I34").Value = 17748.75
$ws.Range("J34").Value = 3311.4092
$ws.Range("K34").Value = 17748.75
$ws.Range("L34").Value = 3311.4092
$ws.Range("M34").Value = -17546.75
$ws.Range("N34").Value = -3715.4092
$ws.Range("H132").Value = 34087.55
$ws.Range("I132").Value = 45308.93
$ws.Range("J132").Value = 7904.3335
$ws.Range("K132").Value = 135926.79
$ws.Range("L132").Value = 23713.0005
$ws.Range("M132").Value = -133396.79
$ws.Range("N132").Value = -28773.0005
$ws.Range("H134").Value = 3809.5588
$ws.Range("I134").Value = 745.4074000000001
$ws.Range("J134").Value = 15628.429
$ws.Range("K134").Value = 2236.2222
$ws.Range("L134").Value = 46885.287
$ws.Range("M134").Value = 298.7777999999998
$ws.Range("N134").Value = -51955.287

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2920
$ws.Range("J80").Value = 2920
$ws.Range("L80").Value = 8760
$ws.Range("N80").Value = -10632
$ws.Range("H83").Value = 2920
$ws.Range("J83").Value = 2920
$ws.Range("L83").Value = 26280
$ws.Range("N83").Value = -35640
$ws.Range("H107").Value = 5074.2085
$ws.Range("I107").Value = 25350
$ws.Range("J107").Value = 1019.05
$ws.Range("K107").Value = 76050
$ws.Range("L107").Value = 3057.15
$ws.Range("M107").Value = -74130
$ws.Range("N107").Value = -6897.15
$ws.Range("H131").Value = 104981.48
$ws.Range("J131").Value = 114456.38
$ws.Range("L131").Value = 343369.14
$ws.Range("N131").Value = -353449.14

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4611.1113
$ws.Range("I113").Value = 4214.2856
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 4214.2856
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -2044.2856
$ws.Range("N113").Value = -10340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7500.6924
$ws.Range("I61").Value = 3750.6667
$ws.Range("J61").Value = 10715
$ws.Range("K61").Value = 3750.6667
$ws.Range("L61").Value = 10715
$ws.Range("M61").Value = -3548.6667
$ws.Range("N61").Value = -11119
$ws.Range("H113").Value = 7500.6924
$ws.Range("I113").Value = 3750.6667
$ws.Range("J113").Value = 10715
$ws.Range("K113").Value = 3750.6667
$ws.Range("L113").Value = 10715
$ws.Range("M113").Value = -1580.6667
$ws.Range("N113").Value = -15055
$ws.Range("H131").Value = 52495
$ws.Range("J131").Value = 52495
$ws.Range("L131").Value = 52495
$ws.Range("N131").Value = -62575
$ws.Range("H132").Value = 1669.0303
$ws.Range("I132").Value = 1132.2222
$ws.Range("K132").Value = 3396.6666
$ws.Range("M132").Value = -866.6665999999996
$ws.Range("H136").Value = 31975.941
$ws.Range("I136").Value = 64435.875
$ws.Range("J136").Value = 3122.6667
$ws.Range("K136").Value = 193307.625
$ws.Range("L136").Value = 9368.000100000001
$ws.Range("M136").Value = -190757.625
$ws.Range("N136").Value = -14468.0001
